# Updated cryptos list on Sat Nov 30 14:15:00 UTC 2024 with GitHub Actions
#
# Refresh the hourly crypto snapshot: new prices + %-change figures for most
# rows, plus a few coin pairs that swapped rank this run (NEAR/PEPE,
# VeChain/Algorand, MantraDAO/WhiteBITCoin).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A leading apostrophe text-qualifies the value so Excel stores it as literal
# text (matching the source inlineStr cells) instead of auto-coercing
# numeric-looking strings such as '241.24' into real numbers.
$ws.Range('D2').Value = '''96.557.65'
$ws.Range('E2').Value = '''  -0.52%  '
$ws.Range('D3').Value = '''3.677.40'
$ws.Range('E3').Value = '''  +1.94%  '
$ws.Range('E4').Value = '''  -0.13%  '
$ws.Range('D5').Value = '''241.24'
$ws.Range('E5').Value = '''  -0.82%  '
$ws.Range('D6').Value = '''1.85'
$ws.Range('E6').Value = '''  +9.61%  '
$ws.Range('D7').Value = '''664.26'
$ws.Range('E7').Value = '''  +0.97%  '
$ws.Range('E8').Value = '''  +1.07%  '
$ws.Range('E9').Value = '''  +1.82%  '
$ws.Range('D10').Value = '''0.999'
$ws.Range('E10').Value = '''  -0.05%  '
$ws.Range('D11').Value = '''3.676.33'
$ws.Range('E11').Value = '''  +1.97%  '
$ws.Range('D12').Value = '''45.69'
$ws.Range('E12').Value = '''  +3.99%  '
$ws.Range('D13').Value = '''0.206'
$ws.Range('E13').Value = '''  +0.54%  '
$ws.Range('D14').Value = '''6.81'
$ws.Range('E14').Value = '''  +5.67%  '
$ws.Range('D15').Value = '''4.362.78'
$ws.Range('E15').Value = '''  +1.99%  '
$ws.Range('E16').Value = '''  +3.14%  '
$ws.Range('D17').Value = '''96.205.53'
$ws.Range('E17').Value = '''  -0.63%  '
$ws.Range('D18').Value = '''8.89'
$ws.Range('E18').Value = '''  +14.52%  '
$ws.Range('D19').Value = '''3.682.00'
$ws.Range('E19').Value = '''  +2.37%  '
$ws.Range('E20').Value = '''  +0.31%  '
$ws.Range('D21').Value = '''18.53'
$ws.Range('E21').Value = '''  +2.67%  '
$ws.Range('D22').Value = '''0.525'
$ws.Range('E22').Value = '''  -1.75%  '
$ws.Range('D23').Value = '''527.38'
$ws.Range('E23').Value = '''  +3.00%  '
$ws.Range('D24').Value = '''3.45'
$ws.Range('E24').Value = '''  +0.83%  '
$ws.Range('B25').Value = '''NEARProtocol'
$ws.Range('C25').Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').Value = '''7.06'
$ws.Range('E25').Value = '''  +2.68%  '
$ws.Range('B26').Value = '''PEPE'
$ws.Range('C26').Value = '''https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').Value = '''0.0000204'
$ws.Range('E26').Value = '''  -0.68%  '
$ws.Range('D27').Value = '''101.92'
$ws.Range('E27').Value = '''  +3.62%  '
$ws.Range('D28').Value = '''13.08'
$ws.Range('E28').Value = '''  -0.22%  '
$ws.Range('D29').Value = '''3.872.57'
$ws.Range('E29').Value = '''  +1.87%  '
$ws.Range('D30').Value = '''0.168'
$ws.Range('E30').Value = '''  +10.24%  '
$ws.Range('D31').Value = '''12.54'
$ws.Range('E31').Value = '''  +6.46%  '
$ws.Range('D32').Value = '''3.05'
$ws.Range('E32').Value = '''  +0.32%  '
$ws.Range('D33').Value = '''1.00'
$ws.Range('E33').Value = '''  +0.17%  '
$ws.Range('D34').Value = '''1.92'
$ws.Range('E34').Value = '''  +18.53%  '
$ws.Range('E35').Value = '''  -0.29%  '
$ws.Range('D36').Value = '''32.76'
$ws.Range('E36').Value = '''  +3.14%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '''  -0.07%  '
$ws.Range('D38').Value = '''651.03'
$ws.Range('E38').Value = '''  +4.63%  '
$ws.Range('D39').Value = '''0.591'
$ws.Range('E39').Value = '''  +3.30%  '
$ws.Range('D40').Value = '''8.86'
$ws.Range('E40').Value = '''  +1.00%  '
$ws.Range('D41').Value = '''44.22'
$ws.Range('E41').Value = '''  +32.19%  '
$ws.Range('E42').Value = '''  +4.65%  '
$ws.Range('D43').Value = '''0.971'
$ws.Range('E43').Value = '''  +4.69%  '
$ws.Range('D44').Value = '''1.98'
$ws.Range('E44').Value = '''  +3.17%  '
$ws.Range('D45').Value = '''6.47'
$ws.Range('E45').Value = '''  +8.32%  '
$ws.Range('E46').Value = '''  +0.01%  '
$ws.Range('B47').Value = '''VeChain'
$ws.Range('C47').Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '''0.0465'
$ws.Range('E47').Value = '''  +6.89%  '
$ws.Range('B48').Value = '''Algorand'
$ws.Range('C48').Value = '''https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').Value = '''0.447'
$ws.Range('E48').Value = '''  +15.35%  '
$ws.Range('D49').Value = '''2.31'
$ws.Range('E49').Value = '''  -0.33%  '
$ws.Range('B50').Value = '''MantraDAO'
$ws.Range('C50').Value = '''https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D50').Value = '''3.71'
$ws.Range('E50').Value = '''  +5.24%  '
$ws.Range('B51').Value = '''WhiteBITCoin'
$ws.Range('C51').Value = '''https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '''23.64'
$ws.Range('E51').Value = '''  -0.20%  '
